# New Microsite scripts support to Beta server
# - Append 3 new interview-history rows (fstcyc167, scndcycle167, finalrun167)
#   to the AMSIN sheet.
# - Append 1 new row (beta167) to the BETA sheet.
# - Fix up the last row (row 34, htfx166tue) on the AMS sheet: tiny timestamp
#   correction + make its cell formatting consistent with the rows above it.

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# AMSIN sheet: add rows 50-52
# ---------------------------------------------------------------------------
$amsin = $wb.Worksheets.Item("AMSIN")

$amsinRows = @(
    @{ Row = 50; A = "2022-09-16"; B = 44820.61133503472; C = "fstcyc167";    D = 155; E = 153; F = 2; G = 3.85 },
    @{ Row = 51; A = "2022-09-19"; B = 44823.62529790509; C = "scndcycle167"; D = 155; E = 153; F = 2; G = 4.81 },
    @{ Row = 52; A = "2022-09-20"; B = 44824.37102759259; C = "finalrun167";  D = 155; E = 154; F = 1; G = 3.84 }
)

foreach ($r in $amsinRows) {
    $row = $r.Row

    # Column A: literal text date (leading apostrophe forces text, same as the
    # existing "Run Date" cells above, instead of Excel auto-converting the
    # string to a date serial)
    $amsin.Range("A$row").Value = "'" + $r.A

    # Column B: numeric timestamp with the sheet's datetime format
    $amsin.Range("B$row").NumberFormat = $dateFmt
    $amsin.Range("B$row").Value = $r.B

    # Column C: sprint name text
    $amsin.Range("C$row").Value = $r.C

    # Columns D-G: numeric counters / duration
    $amsin.Range("D$row").Value = $r.D
    $amsin.Range("E$row").Value = $r.E
    $amsin.Range("F$row").Value = $r.F
    $amsin.Range("G$row").Value = $r.G
}

# ---------------------------------------------------------------------------
# BETA sheet: add row 27
# ---------------------------------------------------------------------------
$beta = $wb.Worksheets.Item("BETA")

$beta.Range("A27").Value = "'2022-09-20"

$beta.Range("B27").NumberFormat = $dateFmt
$beta.Range("B27").Value = 44824.51673351663

$beta.Range("C27").Value = "beta167"

$beta.Range("D27").Value = 155
$beta.Range("E27").Value = 153
$beta.Range("F27").Value = 2
$beta.Range("G27").Value = 3.12

# ---------------------------------------------------------------------------
# AMS sheet: correct row 34 (htfx166tue) timestamp + formatting
# ---------------------------------------------------------------------------
$ams = $wb.Worksheets.Item("AMS")

$ams.Range("A34").Value = "'2022-09-06"

$ams.Range("B34").NumberFormat = $dateFmt
$ams.Range("B34").Value = 44810.92667748842

$ams.Range("C34").Value = "htfx166tue"

$ams.Range("D34").Value = 155
$ams.Range("E34").Value = 148
$ams.Range("F34").Value = 7
$ams.Range("G34").Value = 4.03
